# Weekly update: insert two new daily-price rows for "Ají, Americana (o)"
# at the top of the data block (rows 179-180), pushing the previously
# existing rows 179-193 down to rows 181-195.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the current row 179; this shifts all
# rows from 179 downward by two (old 179 -> 181, ..., old 193 -> 195) and
# keeps Excel's used range/dimension in sync automatically.
$ws.Rows("179:180").Insert()

# Row 179: new record
$ws.Range("A179").Value = 8
$ws.Range("B179").Value = "Terminal La Palmera de La Serena"
$ws.Range("C179").Value = "Coquimbo"
$ws.Range("D179").Value = 44610
$ws.Range("E179").Value = 4
$ws.Range("F179").Value = 100112021
$ws.Range("G179").Value = "Ají"
$ws.Range("H179").Value = "Americana (o)"
$ws.Range("I179").Value = "Primera"
$ws.Range("J179").Value = 580
$ws.Range("K179").Value = 11000
$ws.Range("L179").Value = 12000
$ws.Range("M179").Value = 11500
$ws.Range("N179").Value = '$/caja 15 kilos'
$ws.Range("O179").Value = "Provincia de Limarí"
$ws.Range("P179").Value = 767
$ws.Range("Q179").Value = 15
$ws.Range("R179").Value = "Hortaliza"

# Row 180: new record
$ws.Range("A180").Value = 8
$ws.Range("B180").Value = "Terminal La Palmera de La Serena"
$ws.Range("C180").Value = "Coquimbo"
$ws.Range("D180").Value = 44610
$ws.Range("E180").Value = 4
$ws.Range("F180").Value = 100112021
$ws.Range("G180").Value = "Ají"
$ws.Range("H180").Value = "Americana (o)"
$ws.Range("I180").Value = "Primera"
$ws.Range("J180").Value = 300
$ws.Range("K180").Value = 19000
$ws.Range("L180").Value = 20000
$ws.Range("M180").Value = 19500
$ws.Range("N180").Value = '$/caja 25 kilos'
$ws.Range("O180").Value = "Provincia de Limarí"
$ws.Range("P180").Value = 780
$ws.Range("Q180").Value = 25
$ws.Range("R180").Value = "Hortaliza"
